$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.142.58'
$ws.Range('E2').Value = '  -0.53%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.909.91'
$ws.Range('E3').Value = '  -0.98%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7437'
$ws.Range('E5').Value = '  -1.28%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '245.87'
$ws.Range('E6').Value = '  +0.98%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.9997'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3104'
$ws.Range('E8').Value = '  -2.08%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '26.58'
$ws.Range('E9').Value = '  -5.86%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.06996'
$ws.Range('E10').Value = '  -0.37%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08083'
$ws.Range('E11').Value = '  +0.82%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.7713'
$ws.Range('E12').Value = '  -0.68%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.905.57'
$ws.Range('E13').Value = '  -1.17%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.352'
$ws.Range('E14').Value = '  -0.26%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '92.23'
$ws.Range('E15').Value = '  -1.12%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.42'
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '30.153.06'
$ws.Range('E17').Value = '  -0.54%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '5.988'
$ws.Range('E18').Value = '  +2.63%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007868'
$ws.Range('E19').Value = '  -0.93%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '240.73'
$ws.Range('E20').Value = '  -5.03%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '2.215.75'
$ws.Range('E21').Value = '  +1.59%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '7.170'
$ws.Range('E24').Value = '  +7.08%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.444'
$ws.Range('E25').Value = '  -0.48%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '167.81'
$ws.Range('E26').Value = '  +1.94%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '19.00'
$ws.Range('E27').Value = '  -0.35%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.1294'
$ws.Range('E28').Value = '  -2.87%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.060'
$ws.Range('E29').Value = '  -6.52%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.557'
$ws.Range('E30').Value = '  +2.91%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.350'
$ws.Range('E31').Value = '  -1.13%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.348'
$ws.Range('E32').Value = '  -1.28%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.093'
$ws.Range('E33').Value = '  -0.94%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.319'
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.05168'
$ws.Range('E35').Value = '  -0.26%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7521'
$ws.Range('E36').Value = '  -0.27%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.734'
$ws.Range('E37').Value = '  -2.14%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01953'
$ws.Range('E38').Value = '  -0.13%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.802'
$ws.Range('E39').Value = '  +0.26%  '
$ws.Range('E40').Value = '  -2.47%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.4522'
$ws.Range('E41').Value = '  +0.79%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '74.66'
$ws.Range('E42').Value = '  -4.14%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.997'
$ws.Range('E43').Value = '  +1.29%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.8426'
$ws.Range('E44').Value = '  +1.05%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '7.830'
$ws.Range('E45').Value = '  +3.42%  '
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.001'
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '101.96'
$ws.Range('E47').Value = '  +0.54%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.961'
$ws.Range('E48').Value = '  +0.50%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.109.26'
$ws.Range('E49').Value = '  +1.65%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '36.99'
$ws.Range('E50').Value = '  -1.83%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.1191'
$ws.Range('E51').Value = '  -1.15%  '
